$wb = $excel.ActiveWorkbook

# --- Operator sheet ---
$opSheet = $wb.Worksheets.Item("Operator")
$opSheet.Range("C10").Value = 20
$opSheet.Range("D10").Value = 20

# --- Aircraft sheet ---
$acSheet = $wb.Worksheets.Item("Aircraft")
$acSheet.Range("C2").Value = 200

# --- Ports sheet ---
$portsSheet = $wb.Worksheets.Item("Ports")

# Update sheet view selections / scroll positions
$opSheet.Activate()
$excel.ActiveWindow.ScrollRow = 10
$opSheet.Range("D11").Select()

$acSheet.Activate()
$acSheet.Range("C3").Select()

$portsSheet.Activate()
$excel.ActiveWindow.ScrollRow = 7
$portsSheet.Range("A3").Select()
